$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 301; this pushes the existing rows 301-332
# down to 305-336, matching the target diff exactly.
$ws.Rows("301:304").Insert()

# Copy formatting (incl. the date-cell number format in column D) from the
# now-shifted block (rows 305-308, formerly 301-304) into the freshly
# inserted blank rows so the new block looks identical in style.
$ws.Range("A305:T308").Copy()
$ws.Range("A301").PasteSpecial(-4122)

# Populate the new rows (301-304) with the new weekly data block.
$data = @(
    @(301, 'Especial', 200, 19000, 20000, 19500, '$/caja 10 unidades', 1950, 10),
    @(302, 'Primera',  250, 19000, 20000, 19600, '$/caja 12 unidades', 1633, 12),
    @(303, 'Segunda',  330, 19000, 20000, 19545, '$/caja 14 unidades', 1396, 14),
    @(304, 'Tercera',  300, 19000, 20000, 19500, '$/caja 16 unidades', 1219, 16)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = 'Agrícola del Norte S.A. de Arica'
    $ws.Cells.Item($r, 3).Value = 'Arica y Parinacota'
    $ws.Cells.Item($r, 4).Value = [DateTime]'2023-08-16'
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = 'Fruta'
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = 'Tropicales y subtropicales'
    $ws.Cells.Item($r, 9).Value = 100108005
    $ws.Cells.Item($r, 10).Value = 'Piña'
    $ws.Cells.Item($r, 11).Value = 'Caramelo'
    $ws.Cells.Item($r, 12).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $row[6]
    $ws.Cells.Item($r, 18).Value = 'Ecuador'
    $ws.Cells.Item($r, 19).Value = $row[7]
    $ws.Cells.Item($r, 20).Value = $row[8]
}
